$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Table1")

# 1. Rename "Brown sugar" to "Brown sugar, packed" (both unit-conversion rows)
$ws.Range("B109").Value = "Brown sugar, packed"
$ws.Range("B110").Value = "Brown sugar, packed"

# 2. Add two new ingredient rows to the table: Semolina, Nonfat dried milk
$rowSemolina = $lo.ListRows.Add()
$rSemolina = $rowSemolina.Range
$rSemolina.Cells.Item(1,2).Value = "Semolina"
$rSemolina.Cells.Item(1,3).Value = "Flour"
$rSemolina.Cells.Item(1,5).Value = 1
$rSemolina.Cells.Item(1,6).Value = "cup"
$rSemolina.Cells.Item(1,9).Value = 167
$rSemolina.Cells.Item(1,10).Value = "gram"
$rSemolina.Cells.Item(1,13).Value = "y"
$rSemolina.Cells.Item(1,14).Value = "fdc.nal.usda.gov 10/28/2025"

$rowMilk = $lo.ListRows.Add()
$rMilk = $rowMilk.Range
$rMilk.Cells.Item(1,2).Value = "Nonfat dried milk"
$rMilk.Cells.Item(1,3).Value = "Dairy"
$rMilk.Cells.Item(1,5).Value = 1
$rMilk.Cells.Item(1,6).Value = "cup"
$rMilk.Cells.Item(1,9).Value = 112
$rMilk.Cells.Item(1,10).Value = "gram"
$rMilk.Cells.Item(1,13).Value = "y"
$rMilk.Cells.Item(1,14).Value = "King Arthur 10/28/2025"

# 3. Fill in the ID (UUID) column last, for both new rows
$rSemolina.Cells.Item(1,1).Value = "ae2f976d-b43e-4c3e-8f49-4682a99548ec"
$rMilk.Cells.Item(1,1).Value = "2923bebb-2e9e-46d1-8738-d6109def788a"

# 4. Leave the selection on the newly added cell, matching the end of the edit session
$ws.Range("B133").Select()
